$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the formatting of row 5 down to the new rows 6-16 ---
$src = $ws.Range("A5:K5")
$dst = $ws.Range("A6:K16")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. Clear the existing string cells (H2:K5) so the shared-string table
#        is rebuilt from scratch in the exact order the data is (re)written
#        below (this reproduces the canonical ordering of the shared strings
#        table, matching how the source file was generated). Number/formula
#        cells (A-G) are left untouched since their values do not change
#        type and do not participate in the shared string table.
$ws.Range("H2:K16").ClearContents()

# --- 3. Full row data (table rebuilt for rows 2-16) ---
$rowData = @(
    @{ r = 2;  A = 1;  B = 0; C = 1; D = 0; E = 0; F = 52;  G = 0; H = 'UNKNOWN'; I = 'FAIL';    J = 'Definitely cannot solve, may have got wrong rep.'; K = 'YES' }
    @{ r = 3;  A = 2;  B = 0; C = 1; D = 0; E = 1; F = 21;  G = 0; H = 'SUCCEED'; I = 'FAIL';    J = 'Gets equations right, but solves wrong'; K = 'YES' }
    @{ r = 4;  A = 3;  B = 0; C = 1; D = 1; E = 0; F = 28;  G = 0; H = 'FAIL';    I = 'SUCCEED'; J = 'Gets equations wrong, but it can solve otherwise.'; K = 'YES' }
    @{ r = 5;  A = 4;  B = 0; C = 1; D = 1; E = 1; F = 38;  G = 0; H = 'SUCCEED'; I = 'SUCCEED'; J = 'Non-determinism issues as we do not know why it can''t handle by itself'; K = 'NO' }
    @{ r = 6;  A = 5;  B = 1; C = 1; D = 0; E = 0; F = 82;  G = 0; H = 'UNKNOWN'; I = 'FAIL';    J = 'Non-determinism as we do not know how it could gotten the problem right'; K = 'NO' }
    @{ r = 7;  A = 6;  B = 1; C = 1; D = 0; E = 1; F = 1;   G = 0; H = 'SUCCEED'; I = 'FAIL';    J = 'Can get the right equations, perhaps there is something about the equations from ground truth that mess it up'; K = 'MAYBE' }
    @{ r = 8;  A = 7;  B = 1; C = 1; D = 1; E = 0; F = 45;  G = 0; H = 'SUCCEED'; I = 'SUCCEED'; J = 'ChatGPT is lying on the equations it reports'; K = 'MAYBE' }
    @{ r = 9;  A = 8;  B = 1; C = 1; D = 1; E = 1; F = 5;   G = 0; H = 'SUCCEED'; I = 'SUCCEED'; J = 'Gets everything right'; K = 'YES' }
    @{ r = 10; A = 9;  B = 0; C = 1; D = 0; E = 0; F = 43;  G = 1; H = 'UNKNOWN'; I = 'FAIL';    J = 'Definitely cannot solve, may have got wrong rep.'; K = 'YES' }
    @{ r = 11; A = 10; B = 0; C = 1; D = 0; E = 1; F = 1;   G = 1; H = 'SUCCEED'; I = 'FAIL';    J = 'Gets equations right, but solves wrong'; K = 'YES' }
    @{ r = 12; A = 11; B = 0; C = 1; D = 1; E = 0; F = 18;  G = 1; H = 'FAIL';    I = 'SUCCEED'; J = 'Gets equations wrong, but it can solve otherwise.'; K = 'YES' }
    @{ r = 13; A = 13; B = 1; C = 1; D = 0; E = 0; F = 11;  G = 1; H = 'UNKNOWN'; I = 'FAIL';    J = 'Non-determinism as we do not know how it could gotten the problem right'; K = 'NO' }
    @{ r = 14; A = 14; B = 1; C = 1; D = 0; E = 1; F = 119; G = 1; H = 'SUCCEED'; I = 'FAIL';    J = 'Can get the right equations, perhaps there is something about the equations from ground truth that mess it up'; K = 'MAYBE' }
    @{ r = 15; A = 15; B = 1; C = 1; D = 1; E = 0; F = 27;  G = 1; H = 'SUCCEED'; I = 'SUCCEED'; J = 'ChatGPT is lying on the equations it reports'; K = 'MAYBE' }
    @{ r = 16; A = 16; B = 1; C = 1; D = 1; E = 1; F = 507; G = 1; H = 'SUCCEED'; I = 'SUCCEED'; J = 'Gets everything right'; K = 'YES' }
)

# Columns A-G: plain numbers, safe to (re)write in any order.
foreach ($d in $rowData) {
    $ws.Cells.Item($d.r, 1).Value = $d.A
    $ws.Cells.Item($d.r, 2).Value = $d.B
    $ws.Cells.Item($d.r, 3).Value = $d.C
    $ws.Cells.Item($d.r, 4).Value = $d.D
    $ws.Cells.Item($d.r, 5).Value = $d.E
    $ws.Cells.Item($d.r, 6).Value = $d.F
    $ws.Cells.Item($d.r, 7).Value = $d.G
}

# Columns H and I: only reuse pre-existing text values (UNKNOWN/SUCCEED/FAIL),
# order does not introduce any new shared strings.
foreach ($d in $rowData) {
    $ws.Cells.Item($d.r, 8).Value = $d.H
    $ws.Cells.Item($d.r, 9).Value = $d.I
}

# Column J: must be written top-to-bottom so brand-new text values are
# appended to the shared-string table in the canonical order.
foreach ($d in $rowData) {
    $ws.Cells.Item($d.r, 10).Value = $d.J
}

# Column K: written last, top-to-bottom, so YES/NO/MAYBE are appended after
# all of the column J strings, matching the canonical shared-string order.
foreach ($d in $rowData) {
    $ws.Cells.Item($d.r, 11).Value = $d.K
}

Write-Output "done"
